$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Host ($ws | Get-Member | Select-String "Cells|Range")
